# Updates the cryptos list data (Bitcoin, Ethereum, ... table) on Sheet1.
# For each affected row we may update the Coin name (B), Link (C),
# Price (D) and/or Volume(1h) (E) columns with freshly scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row=2;  B=$null;         C=$null;                                                          D="28.222.93";  E="  +1.89%  " }
    @{ Row=3;  B=$null;         C=$null;                                                          D="1.811.00";   E="  +3.09%  " }
    @{ Row=4;  B=$null;         C=$null;                                                          D="1.002";      E="  -0.37%  " }
    @{ Row=5;  B=$null;         C=$null;                                                          D="324.52";     E="  -0.12%  " }
    @{ Row=6;  B=$null;         C=$null;                                                          D="1.001";      E="  -0.07%  " }
    @{ Row=7;  B=$null;         C=$null;                                                          D="0.4311";     E="  -4.14%  " }
    @{ Row=8;  B=$null;         C=$null;                                                          D="0.3649";     E="  -2.36%  " }
    @{ Row=9;  B=$null;         C=$null;                                                          D="44.73";      E="  -1.30%  " }
    @{ Row=10; B=$null;         C=$null;                                                          D="0.07604";    E="  +0.75%  " }
    @{ Row=11; B=$null;         C=$null;                                                          D="1.136";      E="  +0.88%  " }
    @{ Row=12; B=$null;         C=$null;                                                          D=$null;        E="  -0.01%  " }
    @{ Row=13; B=$null;         C=$null;                                                          D="21.74";      E="  -0.09%  " }
    @{ Row=14; B=$null;         C=$null;                                                          D="6.267";      E="  +0.98%  " }
    @{ Row=15; B=$null;         C=$null;                                                          D="7.423";      E="  +0.99%  " }
    @{ Row=16; B=$null;         C=$null;                                                          D="1.823.99";   E="  +3.86%  " }
    @{ Row=17; B=$null;         C=$null;                                                          D="93.69";      E="  +6.64%  " }
    @{ Row=18; B=$null;         C=$null;                                                          D="0.00001074"; E="  -0.13%  " }
    @{ Row=19; B=$null;         C=$null;                                                          D="0.06393";    E="  +2.80%  " }
    @{ Row=20; B=$null;         C=$null;                                                          D=$null;        E="  -0.05%  " }
    @{ Row=21; B=$null;         C=$null;                                                          D="17.32";      E="  +0.32%  " }
    @{ Row=22; B=$null;         C=$null;                                                          D="6.189";      E="  +0.09%  " }
    @{ Row=23; B=$null;         C=$null;                                                          D="28.231.52";  E="  +1.78%  " }
    @{ Row=24; B=$null;         C=$null;                                                          D="11.51";      E="  -1.36%  " }
    @{ Row=25; B=$null;         C=$null;                                                          D="2.139";      E="  -7.75%  " }
    @{ Row=26; B=$null;         C=$null;                                                          D="160.77";     E="  +5.11%  " }
    @{ Row=27; B=$null;         C=$null;                                                          D="20.55";      E="  -0.53%  " }
    @{ Row=28; B=$null;         C=$null;                                                          D="2.029.37";   E="  +3.83%  " }
    @{ Row=29; B=$null;         C=$null;                                                          D="2.244";      E="  -4.94%  " }
    @{ Row=30; B=$null;         C=$null;                                                          D="130.21";     E="  +1.39%  " }
    @{ Row=31; B=$null;         C=$null;                                                          D="1.182";      E="  -2.92%  " }
    @{ Row=32; B=$null;         C=$null;                                                          D="5.953";      E="  +3.84%  " }
    @{ Row=33; B=$null;         C=$null;                                                          D="0.09089";    E="  -2.50%  " }
    @{ Row=34; B=$null;         C=$null;                                                          D="3.533";      E="  -2.94%  " }
    @{ Row=35; B=$null;         C=$null;                                                          D="12.86";      E="  +1.60%  " }
    @{ Row=36; B=$null;         C=$null;                                                          D=$null;        E="  +2.31%  " }
    @{ Row=37; B=$null;         C=$null;                                                          D="5.169";      E="  +1.78%  " }
    @{ Row=38; B=$null;         C=$null;                                                          D="0.2148";     E="  -0.95%  " }
    @{ Row=39; B=$null;         C=$null;                                                          D="0.6532";     E="  +0.81%  " }
    @{ Row=40; B=$null;         C=$null;                                                          D="0.06156";    E="  +0.14%  " }
    @{ Row=41; B=$null;         C=$null;                                                          D="1.201";      E="  +0.23%  " }
    @{ Row=42; B="WEMIXTOKEN";  C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix";       D="1.425";      E="  +0.32%  " }
    @{ Row=43; B="FraxShare";   C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs";           D="8.001";      E="  +0.04%  " }
    @{ Row=44; B=$null;         C=$null;                                                          D="0.9998";     E="  -0.10%  " }
    @{ Row=45; B=$null;         C=$null;                                                          D="13.80";      E="  -0.02%  " }
    @{ Row=46; B=$null;         C=$null;                                                          D="0.6040";     E="  +0.92%  " }
    @{ Row=47; B=$null;         C=$null;                                                          D="3.726";      E="  -0.72%  " }
    @{ Row=48; B=$null;         C=$null;                                                          D="125.94";     E="  -0.25%  " }
    @{ Row=49; B=$null;         C=$null;                                                          D="2.002";      E="  +0.69%  " }
    @{ Row=50; B=$null;         C=$null;                                                          D="1.172";      E="  +3.59%  " }
    @{ Row=51; B=$null;         C=$null;                                                          D="0.06993";    E="  +1.31%  " }
)

# Force the "Price" column to remain plain text for every row we are about
# to rewrite, so that values such as "28.222.93" or "1.002" are not
# reinterpreted as numbers/dates by Excel when assigned. NumberFormat must
# be applied cell-by-cell; applying it to a multi-area (comma) range only
# affects the first area in this environment.
foreach ($row in $changes) {
    if ($null -ne $row.D) {
        $ws.Range("D" + $row.Row).NumberFormat = "@"
    }
}

foreach ($row in $changes) {
    $r = $row.Row
    if ($null -ne $row.B) {
        $ws.Range("B$r").Value = $row.B
    }
    if ($null -ne $row.C) {
        $ws.Range("C$r").Value = $row.C
    }
    if ($null -ne $row.D) {
        $ws.Range("D$r").Value = $row.D
    }
    if ($null -ne $row.E) {
        $ws.Range("E$r").Value = $row.E
    }
}
